$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 3x3 (sender x target) cluster matrix of NATMI edge-weight stats,
# replacing the previous 6-row (sender x target, excluding FAPs-as-target) table.
$rows = @(
    @("ECs","Plau","St14","ECs",3,1,26.29132066666667,78.87396200000001,0.1411782207947891,0.1411782207947891,3,1,1.051909333333333,3.155728,0.6689916625345733,0.6689916625345734,27.65608559492622,248.904770354336,0.094447052643179,0.09444705264317901),
    @("ECs","Plau","St14","FAPs",3,1,26.29132066666667,78.87396200000001,0.1411782207947891,0.1411782207947891,1,0.3333333333333333,0.1261633333333333,0.37849,0.08023716060215286,0.08023716060215287,3.317000653042222,29.85300587738,0.01132773957543769,0.01132773957543769),
    @("ECs","Plau","St14","sCs",3,1,26.29132066666667,78.87396200000001,0.1411782207947891,0.1411782207947891,3,1,0.3943076666666667,1.182923,0.2507711768632738,0.2507711768632738,10.36686930565845,93.30182375092602,0.03540342857617236,0.03540342857617236),
    @("FAPs","Plau","St14","ECs",3,1,143.4723713333333,430.417114,0.7704129577533824,0.7704129577533824,3,1,1.051909333333333,3.155728,0.6689916625345733,0.6689916625345734,150.9199264809991,1358.279338328992,0.5153998454456132,0.5153998454456133),
    @("FAPs","Plau","St14","FAPs",3,1,143.4723713333333,430.417114,0.7704129577533824,0.7704129577533824,1,0.3333333333333333,0.1261633333333333,0.37849,0.08023716060215286,0.08023716060215287,18.10095260865111,162.90857347786,0.06181574822123775,0.06181574822123775),
    @("FAPs","Plau","St14","sCs",3,1,143.4723713333333,430.417114,0.7704129577533824,0.7704129577533824,3,1,0.3943076666666667,1.182923,0.2507711768632738,0.2507711768632738,56.57225597158022,509.150303744222,0.1931973640865313,0.1931973640865313),
    @("sCs","Plau","St14","ECs",3,1,16.46418733333334,49.39256200000001,0.08840882145182853,0.08840882145182853,3,1,1.051909333333333,3.155728,0.6689916625345733,0.6689916625345734,17.31883232168178,155.869490895136,0.05914476444578102,0.05914476444578103),
    @("sCs","Plau","St14","FAPs",3,1,16.46418733333334,49.39256200000001,0.08840882145182853,0.08840882145182853,1,0.3333333333333333,0.1261633333333333,0.37849,0.08023716060215286,0.08023716060215287,2.077176754597778,18.69459079138,0.007093672805477422,0.007093672805477424),
    @("sCs","Plau","St14","sCs",3,1,16.46418733333334,49.39256200000001,0.08840882145182853,0.08840882145182853,3,1,0.3943076666666667,1.182923,0.2507711768632738,0.2507711768632738,6.491955290969557,58.42759761872601,0.02217038420057009,0.02217038420057009)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $rowdata = $rows[$i]
    for ($j = 0; $j -lt $rowdata.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $rowdata[$j]
    }
}
